$d = $word.ActiveDocument

$titleRange = $d.Paragraphs(1).Range
$titleRange.Find.Execute("Answers: Introduction to radians", $true, $false, $false, $false, $false,
                          $true, 1, $false, "Answers: Introduction to radians", 2)

$authorRange = $d.Paragraphs(2).Range
$authorRange.Find.Execute("Ifan Howells-Baines, Mark Toner", $true, $false, $false, $false, $false,
                           $true, 1, $false, "Ifan Howells-Baines, Mark Toner", 2)

$abstractRange = $d.Paragraphs(4).Range
$abstractRange.Find.Execute("Answers to the questions relating to the guide on radians.", $true, $false, $false, $false, $false,
                             $true, 1, $false, "Answers to the questions relating to the guide on radians.", 2)
